# Atualizado por script em 11-11-2023 20:45
#
# Swaps the two matches recorded on the same "jornada" pairs (the source
# scraper had written them in the wrong order) for rounds played on
# 29/10/2023, 04-05/11/2023 and 09-10/11/2023, and appends the newly
# scraped match Coritiba x Cruzeiro as row 327.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-MatchRow {
    param(
        [int]$Row,
        [string]$Home,
        $HomeGoals,
        [string]$Away,
        $AwayGoals,
        $HomeOpenOdds,
        [string]$HomeOpenDate,
        $HomeCloseOdds,
        [string]$HomeCloseDate,
        $DrawOpenOdds,
        [string]$DrawOpenDate,
        $DrawCloseOdds,
        [string]$DrawCloseDate,
        $AwayOpenOdds,
        [string]$AwayOpenDate,
        $AwayCloseOdds,
        [string]$AwayCloseDate,
        [string]$Url
    )

    $ws.Cells.Item($Row, 6).Value = $Home
    $ws.Cells.Item($Row, 7).Value = $HomeGoals
    $ws.Cells.Item($Row, 8).Value = $Away
    $ws.Cells.Item($Row, 9).Value = $AwayGoals
    $ws.Cells.Item($Row, 10).Value = $HomeOpenOdds
    $ws.Cells.Item($Row, 11).Value = $HomeOpenDate
    $ws.Cells.Item($Row, 12).Value = $HomeCloseOdds
    $ws.Cells.Item($Row, 13).Value = $HomeCloseDate
    $ws.Cells.Item($Row, 14).Value = $DrawOpenOdds
    $ws.Cells.Item($Row, 15).Value = $DrawOpenDate
    $ws.Cells.Item($Row, 16).Value = $DrawCloseOdds
    $ws.Cells.Item($Row, 17).Value = $DrawCloseDate
    $ws.Cells.Item($Row, 18).Value = $AwayOpenOdds
    $ws.Cells.Item($Row, 19).Value = $AwayOpenDate
    $ws.Cells.Item($Row, 20).Value = $AwayCloseOdds
    $ws.Cells.Item($Row, 21).Value = $AwayCloseDate
    $ws.Cells.Item($Row, 22).Value = $Url
}

# --- Row 296 <-> 297 : Corinthians-Santos vs Internacional-Coritiba ---
Set-MatchRow -Row 296 -Home "Internacional" -HomeGoals 3 -Away "Coritiba" -AwayGoals 4 `
    -HomeOpenOdds 1.32 -HomeOpenDate "27/10/2023 02:42" -HomeCloseOdds 1.41 -HomeCloseDate "29/10/2023 22:17" `
    -DrawOpenOdds 5.24 -DrawOpenDate "27/10/2023 02:42" -DrawCloseOdds 4.85 -DrawCloseDate "29/10/2023 22:29" `
    -AwayOpenOdds 10.2 -AwayOpenDate "27/10/2023 02:42" -AwayCloseOdds 8.62 -AwayCloseDate "29/10/2023 22:29" `
    -Url "https://www.betexplorer.com/football/brazil/serie-a/internacional-coritiba/48RKwkFe/"

Set-MatchRow -Row 297 -Home "Corinthians" -HomeGoals 1 -Away "Santos" -AwayGoals 1 `
    -HomeOpenOdds 1.81 -HomeOpenDate "27/10/2023 02:42" -HomeCloseOdds 2.18 -HomeCloseDate "29/10/2023 22:29" `
    -DrawOpenOdds 3.71 -DrawOpenDate "27/10/2023 02:42" -DrawCloseOdds 3.25 -DrawCloseDate "29/10/2023 22:20" `
    -AwayOpenOdds 4.59 -AwayOpenDate "27/10/2023 02:42" -AwayCloseOdds 3.83 -AwayCloseDate "29/10/2023 22:29" `
    -Url "https://www.betexplorer.com/football/brazil/serie-a/corinthians-santos/j1oppixS/"

# --- Row 309 <-> 310 : America MG-Atletico-MG vs Gremio-Bahia ---
Set-MatchRow -Row 309 -Home "Gremio" -HomeGoals 1 -Away "Bahia" -AwayGoals 0 `
    -HomeOpenOdds 1.78 -HomeOpenDate "02/11/2023 00:12" -HomeCloseOdds 1.76 -HomeCloseDate "04/11/2023 23:05" `
    -DrawOpenOdds 3.9 -DrawOpenDate "02/11/2023 00:12" -DrawCloseOdds 3.83 -DrawCloseDate "04/11/2023 23:29" `
    -AwayOpenOdds 4.46 -AwayOpenDate "02/11/2023 00:12" -AwayCloseOdds 4.79 -AwayCloseDate "04/11/2023 23:29" `
    -Url "https://www.betexplorer.com/football/brazil/serie-a/gremio-bahia/QFHwjJ4I/"

Set-MatchRow -Row 310 -Home "America MG" -HomeGoals 1 -Away "Atletico-MG" -AwayGoals 1 `
    -HomeOpenOdds 3.84 -HomeOpenDate "02/11/2023 01:42" -HomeCloseOdds 4.61 -HomeCloseDate "04/11/2023 23:29" `
    -DrawOpenOdds 3.52 -DrawOpenDate "02/11/2023 01:42" -DrawCloseOdds 3.73 -DrawCloseDate "04/11/2023 23:29" `
    -AwayOpenOdds 2.02 -AwayOpenDate "02/11/2023 01:42" -AwayCloseOdds 1.83 -AwayCloseDate "04/11/2023 23:01" `
    -Url "https://www.betexplorer.com/football/brazil/serie-a/america-mg-atletico-mg/rXoh9t5t/"

# --- Row 325 <-> 326 : Botafogo RJ-Gremio vs Bahia-Cuiaba ---
Set-MatchRow -Row 325 -Home "Bahia" -HomeGoals 0 -Away "Cuiaba" -AwayGoals 3 `
    -HomeOpenOdds 1.79 -HomeOpenDate "07/11/2023 01:12" -HomeCloseOdds 1.86 -HomeCloseDate "09/11/2023 23:30" `
    -DrawOpenOdds 3.61 -DrawOpenDate "07/11/2023 01:12" -DrawCloseOdds 3.45 -DrawCloseDate "09/11/2023 23:30" `
    -AwayOpenOdds 4.8 -AwayOpenDate "07/11/2023 01:12" -AwayCloseOdds 4.93 -AwayCloseDate "09/11/2023 23:30" `
    -Url "https://www.betexplorer.com/football/brazil/serie-a/bahia-cuiaba/jcBVMaZu/"

Set-MatchRow -Row 326 -Home "Botafogo RJ" -HomeGoals 3 -Away "Gremio" -AwayGoals 4 `
    -HomeOpenOdds 1.82 -HomeOpenDate "06/11/2023 23:12" -HomeCloseOdds 1.95 -HomeCloseDate "09/11/2023 23:59" `
    -DrawOpenOdds 3.75 -DrawOpenDate "06/11/2023 23:12" -DrawCloseOdds 3.71 -DrawCloseDate "09/11/2023 23:59" `
    -AwayOpenOdds 4.46 -AwayOpenDate "06/11/2023 23:12" -AwayCloseOdds 4.03 -AwayCloseDate "09/11/2023 23:59" `
    -Url "https://www.betexplorer.com/football/brazil/serie-a/botafogo-rj-gremio/Gpp07KZh/"

# --- New row 327 : Coritiba x Cruzeiro ---
# Clone the A/E number formatting (bold+border "Indice" style, date-time
# style) from the row above so the new row matches the rest of the table
# instead of picking up the default "General" style.
$ws.Range("A326").Copy()
$ws.Range("A327").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
$ws.Range("E326").Copy()
$ws.Range("E327").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)

$ws.Cells.Item(327, 1).Value = 326
$ws.Cells.Item(327, 2).Value = "brazil"
$ws.Cells.Item(327, 3).Value = "serie-a"
$ws.Cells.Item(327, 4).Value = "2023"
$ws.Cells.Item(327, 5).Value = 45241.83333333334

Set-MatchRow -Row 327 -Home "Coritiba" -HomeGoals 1 -Away "Cruzeiro" -AwayGoals 0 `
    -HomeOpenOdds 2.99 -HomeOpenDate "08/11/2023 23:12" -HomeCloseOdds 3.63 -HomeCloseDate "11/11/2023 19:42" `
    -DrawOpenOdds 3.28 -DrawOpenDate "08/11/2023 23:12" -DrawCloseOdds 3.35 -DrawCloseDate "11/11/2023 19:42" `
    -AwayOpenOdds 2.51 -AwayOpenDate "08/11/2023 23:12" -AwayCloseOdds 2.2 -AwayCloseDate "11/11/2023 19:42" `
    -Url "https://www.betexplorer.com/football/brazil/serie-a/coritiba-cruzeiro/lE8NAgto/"
